$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new question entry in row 13 (资源管理 / resource log question)
$ws.Cells.Item(13, 2).Value = "资源管理"
$ws.Cells.Item(13, 3).Value = "资源日志的接口没有提供，相关的参数也没有文档。"
$ws.Cells.Item(13, 4).Value = "耿晓红"

# The date column would normally be auto-detected as a real date by Excel;
# enter it as a text formula first and then paste-special as values so it
# lands in the sheet as plain text (matching the other date-like cells
# in this column) rather than a numeric date serial.
$dateCell = $ws.Cells.Item(13, 5)
$dateCell.Formula = '="2015.11.15"'
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)

# The row grew to two lines of text, same as other multi-line rows above.
$ws.Rows("13").RowHeight = 27

# Leave the selection where the author ended up after typing the new row.
$ws.Range("E13").Select()
